# [Kadastro App] Yeni kayit eklendi: 3006
# Adds a new record row (row 59) to both the master "Kayitlar" sheet and
# the filtered "Erdemli" sheet, mirroring the existing layout:
# Kayit No | Tarih | Birim | Parsel Sayisi | Is | Personeller

$wb = $excel.ActiveWorkbook

$newRow = @("3006", "2025-09-11", "Erdemli", "1", "ÇAP", "CEMAL TİMUROĞLU (K.Teknisyeni)")

foreach ($sheetName in @("Kayitlar", "Erdemli")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $targetRange = $ws.Range("A59:F59")

    # Force text formatting first so numeric/date-looking values ("3006",
    # "2025-09-11", "1") are stored as literal text, matching every other
    # cell in the table instead of being auto-converted to a number/date.
    $targetRange.NumberFormat = "@"

    $ws.Range("A59").Value = $newRow[0]
    $ws.Range("B59").Value = $newRow[1]
    $ws.Range("C59").Value = $newRow[2]
    $ws.Range("D59").Value = $newRow[3]
    $ws.Range("E59").Value = $newRow[4]
    $ws.Range("F59").Value = $newRow[5]

    # Drop the temporary number format again so the new cells end up with
    # the same (default) style as all the other data cells.
    $targetRange.ClearFormats()
}

Write-Host "Added record 3006 to Kayitlar and Erdemli sheets"
